$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 450
$ws.Range("I6").Value = 450
$ws.Range("K6").Value = 1350
$ws.Range("M6").Value = -1238

# Row 11
$ws.Range("H11").Value = 15876048
$ws.Range("I11").Value = 15876048
$ws.Range("K11").Value = 15876048
$ws.Range("M11").Value = -15875908

# Row 75
$ws.Range("H75").Value = 350000
$ws.Range("J75").Value = 350000
$ws.Range("L75").Value = 350000
$ws.Range("N75").Value = -351872

# Row 76
$ws.Range("H76").Value = 7939515
$ws.Range("I76").Value = 9262101
$ws.Range("J76").Value = 4000
$ws.Range("K76").Value = 9262101
$ws.Range("L76").Value = 4000
$ws.Range("M76").Value = -9261786
$ws.Range("N76").Value = -4630

# Row 78
$ws.Range("H78").Value = 350000
$ws.Range("J78").Value = 350000
$ws.Range("L78").Value = 1050000
$ws.Range("N78").Value = -1059360

# Row 79
$ws.Range("H79").Value = 7939515
$ws.Range("I79").Value = 9262101
$ws.Range("J79").Value = 4000
$ws.Range("K79").Value = 9262101
$ws.Range("L79").Value = 4000
$ws.Range("M79").Value = -9261009
$ws.Range("N79").Value = -6184

# Row 129
$ws.Range("H129").Value = 1068.8928
$ws.Range("I129").Value = 486.57144
$ws.Range("J129").Value = 1263
$ws.Range("K129").Value = 1459.71432
$ws.Range("L129").Value = 3789
$ws.Range("M129").Value = 3540.28568
$ws.Range("N129").Value = -13789

$ws = $wb.Worksheets.Item("ARM")
# Row 92
$ws.Range("H92").Value = 38000
$ws.Range("J92").Value = 38000
$ws.Range("L92").Value = 38000
$ws.Range("N92").Value = -42992

# Row 102
$ws.Range("H102").Value = 1463.3334
$ws.Range("I102").Value = 1556
$ws.Range("J102").Value = 1000
$ws.Range("K102").Value = 1556
$ws.Range("L102").Value = 1000
$ws.Range("M102").Value = 66
$ws.Range("N102").Value = -4244

# Row 132
$ws.Range("H132").Value = 2117.4468
$ws.Range("I132").Value = 1743.4878
$ws.Range("J132").Value = 4672.8335
$ws.Range("K132").Value = 5230.463400000001
$ws.Range("L132").Value = 14018.5005
$ws.Range("M132").Value = -2700.463400000001
$ws.Range("N132").Value = -19078.5005

# Row 133
$ws.Range("H133").Value = 33874.875
$ws.Range("J133").Value = 33874.875
$ws.Range("L133").Value = 33874.875
$ws.Range("N133").Value = -38934.875

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1243.7037
$ws.Range("I31").Value = 1067.9166
$ws.Range("J31").Value = 2650
$ws.Range("K31").Value = 1067.9166
$ws.Range("L31").Value = 2650
$ws.Range("M31").Value = -772.9166
$ws.Range("N31").Value = -3240

# Row 34
$ws.Range("H34").Value = 1243.7037
$ws.Range("I34").Value = 1067.9166
$ws.Range("J34").Value = 2650
$ws.Range("K34").Value = 1067.9166
$ws.Range("L34").Value = 2650
$ws.Range("M34").Value = -865.9166
$ws.Range("N34").Value = -3054

$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 24447.143
$ws.Range("I4").Value = 94.5
$ws.Range("K4").Value = 283.5
$ws.Range("M4").Value = -171.5

# Row 6
$ws.Range("H6").Value = 250.92308
$ws.Range("I6").Value = 230.16667
$ws.Range("J6").Value = 500
$ws.Range("K6").Value = 690.50001
$ws.Range("L6").Value = 1500
$ws.Range("M6").Value = -577.50001
$ws.Range("N6").Value = -1726

# Row 7
$ws.Range("H7").Value = 166.375
$ws.Range("I7").Value = 110.333336
$ws.Range("J7").Value = 200
$ws.Range("K7").Value = 331.000008
$ws.Range("L7").Value = 600
$ws.Range("M7").Value = -219.000008
$ws.Range("N7").Value = -824

# Row 94
$ws.Range("H94").Value = 2662.5
$ws.Range("J94").Value = 2780
$ws.Range("L94").Value = 8340
$ws.Range("N94").Value = -9692

# Row 98
$ws.Range("H98").Value = 487.5
$ws.Range("I98").Value = 500
$ws.Range("J98").Value = 475
$ws.Range("K98").Value = 1500
$ws.Range("L98").Value = 1425
$ws.Range("M98").Value = -2
$ws.Range("N98").Value = -4421

# Row 103
$ws.Range("H103").Value = 1650.9131
$ws.Range("I103").Value = 637.1
$ws.Range("J103").Value = 2430.7693
$ws.Range("K103").Value = 1911.3
$ws.Range("L103").Value = 7292.3079
$ws.Range("M103").Value = -1032.3
$ws.Range("N103").Value = -9050.3079

# Row 122
$ws.Range("H122").Value = 957.625
$ws.Range("I122").Value = 304
$ws.Range("J122").Value = 1175.5
$ws.Range("K122").Value = 2736
$ws.Range("L122").Value = 10579.5
$ws.Range("M122").Value = -286
$ws.Range("N122").Value = -15479.5

# Row 132
$ws.Range("H132").Value = 5953393
$ws.Range("I132").Value = 756.9091
$ws.Range("J132").Value = 9805099
$ws.Range("K132").Value = 6812.1819
$ws.Range("L132").Value = 88245891
$ws.Range("M132").Value = -4282.1819
$ws.Range("N132").Value = -88250951

# Row 134
$ws.Range("H134").Value = 7526.6
$ws.Range("I134").Value = 5500
$ws.Range("J134").Value = 8877.666999999999
$ws.Range("K134").Value = 16500
$ws.Range("L134").Value = 26633.001
$ws.Range("M134").Value = -11430
$ws.Range("N134").Value = -36773.001

$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 1991.7778
$ws.Range("I102").Value = 1772
$ws.Range("J102").Value = 3750
$ws.Range("K102").Value = 1772
$ws.Range("L102").Value = 3750
$ws.Range("M102").Value = -150
$ws.Range("N102").Value = -6994

# Row 113
$ws.Range("H113").Value = 1624.4546
$ws.Range("I113").Value = 1585.5714
$ws.Range("J113").Value = 1692.5
$ws.Range("K113").Value = 1585.5714
$ws.Range("L113").Value = 1692.5
$ws.Range("M113").Value = 584.4286
$ws.Range("N113").Value = -6032.5

# Row 122
$ws.Range("H122").Value = 742303.9
$ws.Range("I122").Value = 1390189.4
$ws.Range("J122").Value = 1863.2858
$ws.Range("K122").Value = 4170568.2
$ws.Range("L122").Value = 5589.857400000001
$ws.Range("M122").Value = -4168118.2
$ws.Range("N122").Value = -10489.8574

# Row 123
$ws.Range("H123").Value = 9978.956
$ws.Range("I123").Value = 6350
$ws.Range("J123").Value = 10324.571
$ws.Range("K123").Value = 6350
$ws.Range("L123").Value = 10324.571
$ws.Range("M123").Value = -3900
$ws.Range("N123").Value = -15224.571

# Row 132
$ws.Range("H132").Value = 2825.291
$ws.Range("I132").Value = 2431.853
$ws.Range("J132").Value = 3462.2856
$ws.Range("K132").Value = 7295.559
$ws.Range("L132").Value = 10386.8568
$ws.Range("M132").Value = -4765.559
$ws.Range("N132").Value = -15446.8568

$ws = $wb.Worksheets.Item("LTW")
# Row 122
$ws.Range("H122").Value = 3271.0386
$ws.Range("I122").Value = 2468.375
$ws.Range("J122").Value = 3627.7778
$ws.Range("K122").Value = 7405.125
$ws.Range("L122").Value = 10883.3334
$ws.Range("M122").Value = -4955.125
$ws.Range("N122").Value = -15783.3334
